Write-Output ($ppt.Presentations | Get-Member)
